$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# ---------------------------------------------------------------------------
# Shape 2 (id 33, "Rect 0") - box with "23. 그다음 LobbyManager ... 넣어줍니다."
# ---------------------------------------------------------------------------
$shp2 = $s.Shapes.Item(2)

# widen the box slightly (cx 4192270 -> 4192905 EMU)
$shp2.Width = 330.15003937007873

$tr2 = $shp2.TextFrame.TextRange
# merge the "2" + "3" runs into a single "23" run
$tr2.Characters(1,2).Text = "23"
# merge the "그" + "다음 LobbyManager ... 넣어줍니다." runs into one run
$tr2.Characters(5,93).Text = "그다음 LobbyManager 스크립트에 각각의 Button과 InputFiled 그리고 Project 폴더에 Prefab 폴더에 있는 Room 오브젝트를 넣어줍니다."

# ---------------------------------------------------------------------------
# Shape 3 (id 50, "Rect 0") - box with "22. 그리고 룸에 입장했을 때 ..." /
# "마지막으로 룸 목록에 ..."
# ---------------------------------------------------------------------------
$shp3 = $s.Shapes.Item(3)

# widen the box slightly (cx 4104005 -> 4104640 EMU)
$shp3.Width = 323.20003937007874

$tr3 = $shp3.TextFrame.TextRange

# Apply edits from the end of the text backwards so earlier character
# offsets remain valid after each assignment changes the text length.

# merge the final paragraph's two runs ("마지막으로 ... 호출합니다." + " ")
# into a single run with a trailing space
$tr3.Characters(44,81).Text = "마지막으로 룸 목록에 변경 사항이 있을 때 호출되는 함수를 선언하고 룸 생성 함수와 룸 삭제 함수 그리고 룸 목록을 갱신하는 함수를 호출합니다. "

# merge "그리고 " + "룸에 입장했을 때 실패했을 때 호출되는 함수를 선언합니다. "
# into new wording, then re-split into three runs
$tr3.Characters(5,37).Text = "그리고 룸에 입장하는  순간 실패했을 때 호출되는 함수를 선언합니다. "
$tr3.Characters(10,11).Font.Size = 18
$tr3.Characters(21,23).Font.Size = 18

# merge the "2" + "2" runs into a single "22" run
$tr3.Characters(1,2).Text = "22"
